# Adapt column header formatting to the respective input file names:
#   *_old -> *_FV2410   (left / "before" block of the AHB diff)
#   *_new -> *_FV2504   (right / "after" block of the AHB diff)
# Then turn the sheet's data range into a proper Excel Table (with
# autofilter) and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base field names shared by the "_FV2410" (left) and "_FV2504" (right)
# column blocks; column K ("diff") is untouched in between them.
$fieldNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $fieldNames.Count; $i++) {
    $colOld = $i + 1    # columns A..J
    $colNew = $i + 12   # columns L..U
    $ws.Cells.Item(1, $colOld).Value = $fieldNames[$i] + "_FV2410"
    $ws.Cells.Item(1, $colNew).Value = $fieldNames[$i] + "_FV2504"
}

# Turn A1:U94 into an Excel Table (adds the autofilter + tableParts wiring).
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U94"), $null, 1)
$lo.Name = "Table1"

# Freeze the header row (split below row 1, keep the header row visible).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
